$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1549815498154982
$ws.Range("C2").Value = 0.6383763837638377
$ws.Range("J2").Value = 0.01476014760147601
$ws.Range("P2").Value = 0.1143911439114391
$ws.Range("S2").Value = 0.07749077490774908
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.02808988764044944
$ws.Range("J3").Value = 0.03370786516853932
$ws.Range("P3").Value = 0.7134831460674157
$ws.Range("S3").Value = 0.2191011235955056
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.7017543859649122
$ws.Range("S4").Value = 0.2456140350877193
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.01568627450980392
$ws.Range("E6").Value = 0.00392156862745098
$ws.Range("F6").Value = 0.09019607843137255
$ws.Range("J6").Value = 0.2470588235294118
$ws.Range("O6").Value = 0.01568627450980392
$ws.Range("Q6").Value = 0.1764705882352941
$ws.Range("R6").Value = 0.0196078431372549
$ws.Range("S6").Value = 0.3647058823529412
$ws.Range("B7").Value = 0.1473684210526316
$ws.Range("D7").Value = 0.02631578947368421
$ws.Range("F7").Value = 0.05263157894736842
$ws.Range("J7").Value = 0.09473684210526316
$ws.Range("O7").Value = 0.01578947368421053
$ws.Range("Q7").Value = 0.1894736842105263
$ws.Range("R7").Value = 0.05263157894736842
$ws.Range("S7").Value = 0.4210526315789473
$ws.Range("B8").Value = 0.08979591836734693
$ws.Range("D8").Value = 0.03061224489795918
$ws.Range("F8").Value = 0.08571428571428572
$ws.Range("J8").Value = 0.06938775510204082
$ws.Range("O8").Value = 0.01224489795918367
$ws.Range("Q8").Value = 0.1959183673469388
$ws.Range("R8").Value = 0.09591836734693877
$ws.Range("S8").Value = 0.4204081632653061
$ws.Range("B9").Value = 0.08648648648648649
$ws.Range("D9").Value = 0.02162162162162162
$ws.Range("F9").Value = 0.07567567567567568
$ws.Range("J9").Value = 0.1081081081081081
$ws.Range("O9").Value = 0.01621621621621622
$ws.Range("Q9").Value = 0.2702702702702703
$ws.Range("R9").Value = 0.06486486486486487
$ws.Range("S9").Value = 0.3567567567567568
$ws.Range("B10").Value = 0.1012765957446809
$ws.Range("D10").Value = 0.02638297872340425
$ws.Range("E10").Value = 0.001702127659574468
$ws.Range("F10").Value = 0.06808510638297872
$ws.Range("J10").Value = 0.1046808510638298
$ws.Range("O10").Value = 0.01446808510638298
$ws.Range("Q10").Value = 0.2459574468085106
$ws.Range("R10").Value = 0.06978723404255319
$ws.Range("S10").Value = 0.3676595744680851
$ws.Range("G11").Value = 0.1455696202531646
$ws.Range("J11").Value = 0.1265822784810127
$ws.Range("K11").Value = 0.2025316455696203
$ws.Range("L11").Value = 0.5158227848101266
$ws.Range("S11").Value = 0.00949367088607595
$ws.Range("G12").Value = 0.7484276729559748
$ws.Range("J12").Value = 0.2138364779874214
$ws.Range("K12").Value = 0.006289308176100629
$ws.Range("L12").Value = 0.006289308176100629
$ws.Range("S12").Value = 0.02515723270440252
$ws.Range("F13").Value = 0.02083333333333333
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.007782101167315175
$ws.Range("H15").Value = 0.1673151750972763
$ws.Range("I15").Value = 0.09727626459143969
$ws.Range("J15").Value = 0.3424124513618677
$ws.Range("K15").Value = 0.0622568093385214
$ws.Range("M15").Value = 0.01167315175097276
$ws.Range("O15").Value = 0.05836575875486381
$ws.Range("S15").Value = 0.2529182879377432
$ws.Range("F16").Value = 0.03061224489795918
$ws.Range("H16").Value = 0.173469387755102
$ws.Range("I16").Value = 0.02551020408163265
$ws.Range("J16").Value = 0.4336734693877551
$ws.Range("K16").Value = 0.1581632653061225
$ws.Range("M16").Value = 0.01530612244897959
$ws.Range("O16").Value = 0.04591836734693878
$ws.Range("S16").Value = 0.1173469387755102
$ws.Range("F17").Value = 0.02348336594911937
$ws.Range("H17").Value = 0.2093933463796477
$ws.Range("I17").Value = 0.08806262230919765
$ws.Range("J17").Value = 0.3913894324853229
$ws.Range("K17").Value = 0.08610567514677103
$ws.Range("M17").Value = 0.01956947162426614
$ws.Range("O17").Value = 0.08023483365949119
$ws.Range("S17").Value = 0.1017612524461839
$ws.Range("F18").Value = 0.02580645161290323
$ws.Range("H18").Value = 0.2064516129032258
$ws.Range("I18").Value = 0.0967741935483871
$ws.Range("J18").Value = 0.3290322580645161
$ws.Range("K18").Value = 0.1161290322580645
$ws.Range("M18").Value = 0.05161290322580645
$ws.Range("O18").Value = 0.07096774193548387
$ws.Range("S18").Value = 0.1032258064516129
$ws.Range("F19").Value = 0.0288
$ws.Range("H19").Value = 0.2208
$ws.Range("I19").Value = 0.0768
$ws.Range("J19").Value = 0.332
$ws.Range("K19").Value = 0.1136
$ws.Range("M19").Value = 0.0208
$ws.Range("O19").Value = 0.092
$ws.Range("S19").Value = 0.1152
